# Pith2Bark/dataset/tracy.xlsx -- "add classification and update structure"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Row -> [B value, C value]
$data = @{
    2  = @(124, 130)
    3  = @(120, 109)
    4  = @(244, 222)
    5  = @(130, 116)
    6  = @(125, 113)
    7  = @(143, 124)
    8  = @(129, 118)
    9  = @(114, 124)
    10 = @(90,  80)
    11 = @(79,  68)
    12 = @(107, 96)
    13 = @(90,  98)
    14 = @(83,  72)
    15 = @(84,  79)
    16 = @(83,  80)
    17 = @(89,  77)
    18 = @(88,  79)
    19 = @(78,  69)
    20 = @(64,  59)
    21 = @(63,  56)
    22 = @(80,  92)
    23 = @(101, 94)
}

# Rows 12-23 don't have a formatted "Count 2" (column C) cell yet, so bring
# the column C formatting in line with column B (which already carries the
# sheet's data style) before filling the new classification counts in.
$ws.Range("B12:B23").Copy()
$ws.Range("C12:C23").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
}

# Move/save the selection at D23 (matches <selection activeCell="D23" sqref="D23"/>)
$ws.Range("D23").Select()
